$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
$ws.Range("A1").Value = "Metrica"
$ws.Range("B1").Value = "N_Manual"
$ws.Range("C1").Value = "N_IA"
$ws.Range("D1").Value = "Media_Manual"
$ws.Range("E1").Value = "Media_IA"
$ws.Range("F1").Value = "SD_Manual"
$ws.Range("G1").Value = "SD_IA"
$ws.Range("H1").Value = "Diferencia"
$ws.Range("I1").Value = "t_statistic_std"
$ws.Range("J1").Value = "p_value_std"
$ws.Range("K1").Value = "t_statistic_welch"
$ws.Range("L1").Value = "p_value_welch"
$ws.Range("M1").Value = "Cohens_d"
$ws.Range("N1").Value = "Interpretacion"
$ws.Range("O1").Value = "Significativo"

# --- Row 2: Instruction Coverage (%) ---
$ws.Range("A2").Value = "Instruction Coverage (%)"
$ws.Range("B2").Value = 6
$ws.Range("C2").Value = 6
$ws.Range("D2").Value = 18.25166666666667
$ws.Range("E2").Value = 17.67525
$ws.Range("F2").Value = 12.50267717997496
$ws.Range("G2").Value = 11.36214431236463
$ws.Range("H2").Value = 0.5764166666666632
$ws.Range("I2").Value = 0.08357445633243105
$ws.Range("J2").Value = 0.9350441276338065
$ws.Range("K2").Value = 0.08357445633243103
$ws.Range("L2").Value = 0.9350588898977312
$ws.Range("M2").Value = 0.04825173486090568
$ws.Range("N2").Value = "Negligible"
$ws.Range("O2").Value = "NO"

# --- Row 3: Branch Coverage (%) ---
$ws.Range("A3").Value = "Branch Coverage (%)"
$ws.Range("B3").Value = 6
$ws.Range("C3").Value = 6
$ws.Range("D3").Value = 14.58333333333333
$ws.Range("E3").Value = 12.046875
$ws.Range("F3").Value = 12.61612724518371
$ws.Range("G3").Value = 9.491644145971234
$ws.Range("H3").Value = 2.536458333333334
$ws.Range("I3").Value = 0.3935308886453862
$ws.Range("J3").Value = 0.702187336897783
$ws.Range("K3").Value = 0.3935308886453862
$ws.Range("L3").Value = 0.7028099954649281
$ws.Range("M3").Value = 0.2272051644938464
$ws.Range("N3").Value = "Pequeño"
$ws.Range("O3").Value = "NO"

# --- Row 4: Mutation Score (%) ---
$ws.Range("A4").Value = "Mutation Score (%)"
$ws.Range("B4").Value = 6
$ws.Range("C4").Value = 6
$ws.Range("D4").Value = 18.52
$ws.Range("E4").Value = 14.757125
$ws.Range("F4").Value = 17.71322782555455
$ws.Range("G4").Value = 11.61356661813028
$ws.Range("H4").Value = 3.762874999999999
$ws.Range("I4").Value = 0.4351606632835607
$ws.Range("J4").Value = 0.6726876222917011
$ws.Range("K4").Value = 0.4351606632835607
$ws.Range("L4").Value = 0.6741269956881768
$ws.Range("M4").Value = 0.2512401260874998
$ws.Range("N4").Value = "Pequeño"
$ws.Range("O4").Value = "NO"

# --- Row 5: Time (seconds) ---
$ws.Range("A5").Value = "Time (seconds)"
$ws.Range("B5").Value = 6
$ws.Range("C5").Value = 6
$ws.Range("D5").Value = 0.0822592261904762
$ws.Range("E5").Value = 0.1937777083333333
$ws.Range("F5").Value = 0.06882194001827036
$ws.Range("G5").Value = 0.1931202829042315
$ws.Range("H5").Value = -0.1115184821428571
$ws.Range("I5").Value = -1.332394857396161
$ws.Range("J5").Value = 0.2122952491574715
$ws.Range("K5").Value = -1.332394857396161
$ws.Range("L5").Value = 0.2292531524231591
$ws.Range("M5").Value = -0.7692585295845469
$ws.Range("N5").Value = "Mediano"
$ws.Range("O5").Value = "NO"

# --- Remove old column P content (was "direccion" column, now removed) ---
$ws.Range("P1:P5").Clear()
